$d = $word.ActiveDocument

# Remove all tables in the document
while ($d.Tables.Count -gt 0) {
    $d.Tables.Item(1).Delete()
}

# Remove all remaining paragraph content (text + paragraph marks except
# the final mandatory paragraph mark Word always keeps)
$r = $d.Range(0, $d.Content.End)
$r.Delete()
